$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking text (e.g. "0.530", "10.00")
# must be forced to Text so Excel keeps the exact string instead of
# coercing/parsing them into a number (which would drop trailing zeros, etc).
$textCells = @("D5","D6","D9","D12","D13","D20","D21","D22","D23","D24","D25","D27","D28","D29","D30","D31","D33","D39","D40","D41","D43","D44","D47","D48","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "68.643.99"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "3.847.69"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "600.71"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Value = "163.72"
$ws.Range("E6").Value = "  -2.20%  "
$ws.Range("D7").Value = "3.846.14"
$ws.Range("E7").Value = "  +2.86%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "36.94"
$ws.Range("E13").Value = "  -3.12%  "
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").Value = "4.490.01"
$ws.Range("E15").Value = "  +2.84%  "
$ws.Range("D16").Value = "3.864.29"
$ws.Range("E16").Value = "  +3.38%  "
$ws.Range("D17").Value = "68.801.93"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("E18").Value = "  +2.69%  "
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").Value = "17.12"
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").Value = "11.28"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "485.44"
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("D23").Value = "0.719"
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").Value = "0.0000160"
$ws.Range("E24").Value = "  +6.51%  "
$ws.Range("D25").Value = "84.06"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("D27").Value = "12.14"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "10.00"
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "2.95"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "7.90"
$ws.Range("E31").Value = "  -4.09%  "
$ws.Range("B32").Value = "WrappedeETH"
$ws.Range("C32").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D32").Value = "3.996.52"
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("D33").Value = "2.38"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("E34").Value = "  +1.61%  "
$ws.Range("D35").Value = "3.791.70"
$ws.Range("E35").Value = "  +3.22%  "
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").Value = "5.88"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "0.319"
$ws.Range("E41").Value = "  -2.34%  "
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("D43").Value = "432.82"
$ws.Range("E43").Value = "  +2.21%  "
$ws.Range("D44").Value = "48.51"
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D47").Value = "8.40"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("D48").Value = "142.92"
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("D49").Value = "2.840.13"
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0357"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "25.77"
$ws.Range("E51").Value = "  +13.24%  "
